# Rerun and summarise models without urban landuse:
#  - rename each "summ*" sheet to its new id
#  - update the "Education[T.Unknown]" label to "Education[T.Unknown/Other]"
#    in row 5 of every sheet

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ40133527",
    "summ40863691",
    "summ41665455",
    "summ42415121",
    "summ43173716",
    "summ43948294",
    "summ44722097",
    "summ45492432",
    "summ46281449"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]

    if ($ws.Range("A5").Text -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
